$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.611.58'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '3.229.03'
$ws.Range('E3').Value = '  +1.10%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.42'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.225.80'
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.70'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -5.66%  '
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000272'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '39.03'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = '3.763.20'
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '66.696.96'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '3.232.36'
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('E19').Value = '  +1.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '511.06'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.03'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.65'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.06'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.01'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.09'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.38'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +4.13%  '
$ws.Range('E30').Value = '  +2.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.03'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.15'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.113'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +26.94%  '
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('E35').Value = '  -3.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.49'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '507.27'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.45'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').Value = '0.0₃0772'
$ws.Range('E39').Value = '  +16.87%  '
$ws.Range('E40').Value = '  +0.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.05'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +7.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.130'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.53%  '
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('E45').Value = '  +1.79%  '
$ws.Range('D46').Value = '2.891.16'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.35'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('E48').Value = '  +4.14%  '
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '122.81'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.65%  '
